$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes the old rows 4-13 down to 5-14).
# Excel adjusts all relative formulas automatically (e.g. G12's
# "=G2+G4+G6+G9" becomes "=G2+G5+G7+G10", matching the commit diff).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the "BOSE inputs" data block.
$ws.Cells.Item(4, 1).Value = "BOSE inputs"

$ws.Cells.Item(4, 2).Value = 4
$ws.Cells.Item(4, 2).ClearFormats()

$ws.Cells.Item(4, 3).Formula = "=PI()*(0.005/2)^2*25.4^2"

$ws.Cells.Item(4, 4).Value = 1

$ws.Cells.Item(4, 5).Formula = "=(C4/1000^2)/D4*1000"

$ws.Cells.Item(4, 6).Formula = "=10000/1000"
$ws.Cells.Item(4, 6).ClearFormats()

$ws.Cells.Item(4, 7).Formula = "=B4*E4*F4*1000000"

$ws.Cells.Item(4, 8).Value = "36 AWG Phosphor Bronze - see Lakeshore appendix for data value"

# New data further down the sheet (feedhorn-drawing related numbers).
$ws.Cells.Item(17, 6).Value = 49.9818
$ws.Cells.Item(17, 7).Value = 0.101481

$ws.Cells.Item(19, 6).Formula = "=F17/G17"

# Author finished with the cursor on the new G4 cell.
$ws.Range("G4").Select() | Out-Null
